$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.075.26'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '2.941.13'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '''354.96'
$ws.Range("E5").Value = '  +0.93%  '
$ws.Range("D6").Value = '''108.20'
$ws.Range("E6").Value = '  -4.29%  '
$ws.Range("D7").Value = '''0.565'
$ws.Range("E7").Value = '  +1.74%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '''0.622'
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").Value = '''38.40'
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("D12").Value = '''0.0866'
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").Value = '''19.27'
$ws.Range("E13").Value = '  -2.52%  '
$ws.Range("D14").Value = '3.423.91'
$ws.Range("E14").Value = '  +1.91%  '
$ws.Range("D15").Value = '''7.71'
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("D16").Value = '2.928.30'
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("D17").Value = '''0.976'
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("D18").Value = '52.029.13'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").Value = '''3.47'
$ws.Range("E19").Value = '  +4.31%  '
$ws.Range("D20").Value = '''7.54'
$ws.Range("E20").Value = '  -0.70%  '
$ws.Range("D21").Value = '''13.67'
$ws.Range("E21").Value = '  -1.56%  '
$ws.Range("D22").Value = '0.0₃0976'
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").Value = '''70.01'
$ws.Range("E23").Value = '  -1.61%  '
$ws.Range("D24").Value = '''266.42'
$ws.Range("E24").Value = '  -1.09%  '
$ws.Range("D25").Value = '''2.76'
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("D26").Value = '''0.177'
$ws.Range("E26").Value = '  -2.72%  '
$ws.Range("D27").Value = '''26.95'
$ws.Range("E27").Value = '  +0.77%  '
$ws.Range("D28").Value = '''7.63'
$ws.Range("E28").Value = '  +14.81%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("D31").Value = '''10.35'
$ws.Range("E31").Value = '  -2.72%  '
$ws.Range("D32").Value = '''36.67'
$ws.Range("E32").Value = '  -2.11%  '
$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D33").Value = '''6.03'
$ws.Range("E33").Value = '  -3.02%  '
$ws.Range("B34").Value = 'Toncoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D34").Value = '''2.18'
$ws.Range("E34").Value = '  -3.76%  '
$ws.Range("D35").Value = '''52.12'
$ws.Range("E35").Value = '  -1.81%  '
$ws.Range("D36").Value = '''0.0436'
$ws.Range("E36").Value = '  -2.87%  '
$ws.Range("D37").Value = '''0.997'
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").Value = '''3.18'
$ws.Range("E38").Value = '  -3.49%  '
$ws.Range("E39").Value = '  -1.88%  '
$ws.Range("D40").Value = '''17.98'
$ws.Range("E40").Value = '  -4.33%  '
$ws.Range("D41").Value = '''2.70'
$ws.Range("E41").Value = '  -0.97%  '
$ws.Range("D42").Value = '''0.119'
$ws.Range("E42").Value = '  +1.55%  '
$ws.Range("D43").Value = '''23.36'
$ws.Range("E43").Value = '  +1.78%  '
$ws.Range("D44").Value = '''118.62'
$ws.Range("E44").Value = '  -0.35%  '
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("D46").Value = '''2.47'
$ws.Range("E46").Value = '  -3.49%  '
$ws.Range("D47").Value = '2.123.29'
$ws.Range("E47").Value = '  -2.19%  '
$ws.Range("D48").Value = '''3.38'
$ws.Range("E48").Value = '  -3.00%  '
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '3.255.41'
$ws.Range("E49").Value = '  +2.00%  '
$ws.Range("D50").Value = '''0.242'
$ws.Range("E50").Value = '  -7.97%  '
$ws.Range("B51").Value = 'BEAM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D51").Value = '''0.0346'
$ws.Range("E51").Value = '  +0.59%  '
